$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2 ---
$ws.Cells.Item(2, 56).Value = 151    # BD2: 126 -> 151

# --- Row 5 ---
$ws.Cells.Item(5, 7).Value = 8       # G5:  8.5  -> 8
$ws.Cells.Item(5, 8).Value = 3.9     # H5:  4    -> 3.9
$ws.Cells.Item(5, 9).Value = 1.48    # I5:  1.45 -> 1.48
$ws.Cells.Item(5, 10).Value = 7.5    # J5:  8    -> 7.5
$ws.Cells.Item(5, 25).Value = 23     # Y5:  26   -> 23
$ws.Cells.Item(5, 30).Value = 7.5    # AD5: 8    -> 7.5
$ws.Cells.Item(5, 36).Value = 9.5    # AJ5: 9    -> 9.5

# --- Row 6 ---
$ws.Cells.Item(6, 12).Value = 10     # L6:  11 -> 10
$ws.Cells.Item(6, 30).Value = 11     # AD6: 10 -> 11

# --- Row 7 ---
$ws.Cells.Item(7, 7).Value = 4.1     # G7:  4    -> 4.1
$ws.Cells.Item(7, 8).Value = 3.9     # H7:  3.8  -> 3.9
$ws.Cells.Item(7, 9).Value = 1.75    # I7:  1.8  -> 1.75
$ws.Cells.Item(7, 10).Value = 4.5    # J7:  4.33 -> 4.5
$ws.Cells.Item(7, 12).Value = 2.38   # L7:  2.4  -> 2.38
$ws.Cells.Item(7, 21).Value = 1.73   # U7:  1.67 -> 1.73
$ws.Cells.Item(7, 22).Value = 2      # V7:  2.1  -> 2
$ws.Cells.Item(7, 33).Value = 8      # AG7: 8.5  -> 8
$ws.Cells.Item(7, 34).Value = 9      # AH7: 9.5  -> 9
$ws.Cells.Item(7, 42).Value = 29     # AP7: 26   -> 29
$ws.Cells.Item(7, 47).Value = 8      # AU7: 7.5  -> 8
$ws.Cells.Item(7, 50).Value = 9      # AX7: 9.5  -> 9

# --- Remove the last three match rows (14-16) entirely ---
# This also shrinks the used range from A1:BD16 down to A1:BD13.
$ws.Rows("14:16").Delete()
